# Update "想去人数" (number of people wanting to attend) values that were
# refreshed by the data-scraping job (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 2 and 3
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 976
$wsExpo.Range("F3").Value = 1973

# Sheet "全部类型" (all categories) - rows 4 and 5 mirror the same events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 976
$wsAll.Range("F5").Value = 1973
